$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update C184: "Reconsile Revision" -> "Reconsile Revision & changes" ---
$ws.Range("C184").Value = "Reconsile Revision & changes"

# --- Add a new daily block (rows 189-194), mirroring the 181-187 block layout ---
# Copy formatting from the previous block (181:187) down to the new block start (189),
# which also creates a trailing extra row (195) that gets removed afterwards.
$ws.Range("A181:D187").Copy()
$ws.Range("A189").PasteSpecial(-4122)
$ws.Rows("195").Delete()

# Row 189: date + "Domm"
$ws.Range("A189").Value = 45700
$ws.Range("B189").Value = "Domm"
$ws.Range("D189").Value = 0.25

# Row 190: "Meeting" / "Reconsile"
$ws.Range("B190").Value = "Meeting"
$ws.Range("C190").Value = "Reconsile"
$ws.Range("D190").Value = 1

# Row 191: "General Discussion"
$ws.Range("C191").Value = "General Discussion"
$ws.Range("D191").Value = 0.25

# Row 192: "Study" / "Reconsile Revision & changes"
$ws.Range("B192").Value = "Study"
$ws.Range("C192").Value = "Reconsile Revision & changes"
$ws.Range("D192").Value = 2

# Row 193: "Editors – NumberBox, SelectBox"
$ws.Range("C193").Value = "Editors – NumberBox, SelectBox"
$ws.Range("D193").Value = 4.5

# Row 194: Total
$ws.Range("B194").Value = "Total"
$ws.Range("D194").Formula = "=SUM(D188:D193)"
